$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the subcategory (column H) text for the rows whose label changed.
$ws.Range("H2").Value = "data collection, data analysis, data gathering diagram"
$ws.Range("H3").Value = "data collection, data analysis, data gathering diagram"
$ws.Range("H4").Value = "data collection, data analysis, data gathering diagram"
$ws.Range("H5").Value = "data collection, data analysis, data gathering diagram"
$ws.Range("H6").Value = "data collection, data analysis, data gathering diagram"
$ws.Range("H7").Value = "data collection, data analysis, data gathering diagram"
$ws.Range("H8").Value = "drawing(s)"
$ws.Range("H9").Value = "drawing(s)"
$ws.Range("H12").Value = "line graph(s)"
$ws.Range("H13").Value = "scatter plot(s)"
$ws.Range("H14").Value = "scatter plot(s)"
$ws.Range("H15").Value = "line graph(s)"
$ws.Range("H16").Value = "line graph(s)"
$ws.Range("H17").Value = "line graph(s)"
$ws.Range("H18").Value = "data display"
$ws.Range("H20").Value = "line graph(s)"
$ws.Range("H30").Value = "line graph(s)"

# Remove the entire "is_viewed" column (column I), which shifts nothing else
# and tightens the sheet's dimension down to A1:H35.
$ws.Columns("I:I").Delete()
